$d = $word.ActiveDocument

$replacements = @(
    @("Resumen de transformación digital", "Resumen de la transformación digital"),
    @("Fabrikam Inc. ha estado experimentando una iniciativa de transformación digital integral destinada a mejorar la eficiencia operativa, mejorar la experiencia del cliente y impulsar la innovación.", "Fabrikam Inc. ha estado inmersa en una iniciativa de transformación digital integral destinada a mejorar la eficiencia operativa, mejorar la experiencia del cliente e impulsar la innovación."),
    @("En el resumen siguiente se describen las actualizaciones clave y los hitos alcanzados hasta la fecha.", "En el siguiente resumen se describen las actualizaciones clave y los hitos alcanzados hasta la fecha."),
    @("Actualizaciones de claves", "Actualizaciones claves"),
    @("Implementación de la infraestructura en la nube", "Implementación de una infraestructura en la nube"),
    @("Migración del 80 % de las aplicaciones locales a la nube.", "Se ha migrado el 80 % de las aplicaciones locales a la nube."),
    @("Escalabilidad mejorada y reducción de los costos de TI en un 25 %.", "Se ha mejorado la escalabilidad y se han reducido los costes de TI en un 25 %."),
    @("Seguridad y cumplimiento mejorados de los datos con los estándares del sector.", "Se han mejorado el cumplimiento y la seguridad de los datos con los estándares del sector."),
    @("Adopción de la inteligencia artificial y el aprendizaje automático", "Adopción de la IA y el aprendizaje automático"),
    @("Integración del análisis con tecnología de inteligencia artificial para simplificar los procesos de toma de decisiones.", "Se ha integrado el análisis con tecnología de IA para optimizar los procesos de toma de decisiones."),
    @("Implementación de modelos de Machine Learning para predecir el comportamiento del cliente y personalizar los esfuerzos de marketing.", "Se han implementado modelos de aprendizaje automático para predecir el comportamiento del cliente y personalizar los esfuerzos de marketing."),
    @("Reducción de los procesos manuales, lo que conduce a un aumento del 30 % en la productividad.", "Se han reducido los procesos manuales, lo que ha llevado a un aumento del 30 % en la productividad."),
    @("Inicio de un nuevo portal de clientes con funcionalidades de autoservicio.", "Se ha lanzado un nuevo portal de clientes con funcionalidades de autoservicio."),
    @("Introducción de bots de chat para el soporte al cliente del 24/7, lo que reduce los tiempos de respuesta en un 50 %.", "Se han introducido bots de chat para brindar asistencia al cliente 24 horas al día, 7 días a la semana, lo que ha llevado a reducir los tiempos de respuesta en un 50 %."),
    @("Se mejoraron las puntuaciones de satisfacción del cliente en un 20 % durante el último año.", "Se han mejorado las puntuaciones de satisfacción del cliente en un 20 % durante el último año."),
    @("Implementación de la automatización de procesos robóticos (RPA) para tareas rutinarias.", "Se ha implementado la automatización robótica de procesos (RPA) para tareas rutinarias."),
    @("Reasignación de recursos humanos a roles más estratégicos dentro de la organización.", "Se han reasignado los recursos humanos a roles más estratégicos dentro de la organización."),
    @("Formación y desarrollo de empleados", "Formación y desarrollo de los empleados"),
    @("Llevó a cabo programas de alfabetización digital para todos los empleados.", "Se han llevado a cabo programas de alfabetización digital para todos los empleados."),
    @("Lanzó una nueva plataforma de aprendizaje electrónico con cursos sobre tecnologías emergentes.", "Se ha lanzado una nueva plataforma de aprendizaje electrónico con cursos sobre tecnologías emergentes."),
    @("Aumento del compromiso de los empleados y la adopción de nuevas herramientas en un 35 %.", "Se ha aumentado el compromiso del empleado y la adopción de nuevas herramientas en un 35 %."),
    @("Q1 2024: Migración completada a la infraestructura en la nube.", "T1 2024: Se completó la migración a la infraestructura en la nube."),
    @("Q2 2024: Lanzamiento de la plataforma de análisis controlada por IA.", "T2 2024: Lanzamiento de una plataforma de análisis impulsada por IA."),
    @("Q3 2024: Se introdujo un nuevo portal de clientes digitales.", "T3 2024: Se introdujo un nuevo portal del cliente digital."),
    @("Q4 2024: Se ha logrado la automatización del 50 % de los procesos rutinarios.", "T4 2024: Se logró la automatización del 50 % de los procesos rutinarios."),
    @("Continúe expandiendo las aplicaciones de inteligencia artificial y aprendizaje automático en todos los departamentos.", "Continuar expandiendo las aplicaciones de IA y aprendizaje automático en todos los departamentos."),
    @("Mejore aún más la experiencia del cliente digital con nuevas características y servicios.", "Mejorar aún más la experiencia del cliente digital con nuevas características y servicios."),
    @("Céntrese en las medidas de ciberseguridad para protegerse frente a amenazas en evolución.", "Centrar la atención en las medidas de ciberseguridad para tener protección frente a amenazas en constante evolución."),
    @("Desarrolle una estrategia digital completa para los próximos cinco años.", "Desarrollar una estrategia digital completa para los próximos cinco años."),
    @("La organización se compromete a aprovechar la tecnología para impulsar el crecimiento futuro y la innovación.", "La organización sigue comprometida a aprovechar la tecnología para impulsar el crecimiento futuro y la innovación."),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

# Bold fixes for two Heading runs that should be bold
$boldHeadings = @("Automatización de procesos", "Planes de futuro")
foreach ($headingText in $boldHeadings) {
    $range = $d.Content
    $found = $range.Find.Execute($headingText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = 1
    } else {
        Write-Host "HEADING NOT FOUND: $headingText"
    }
}
